# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with freshly scraped values (GitHub Actions cron refresh).
#
# D-column values that look numeric (e.g. "1.006") get auto-coerced to a
# real number by Excel's normal Value-assignment heuristics, which would
# change their stored cell type from text to numeric. The source data is
# text (note the non-numeric entries like "25.881.53" which use dots as
# thousands separators), so for any new D value that parses as a plain
# number we briefly force the cell to Text format, assign the value, and
# restore the default "Normal" style/format afterwards so the cell keeps
# looking exactly like its untouched neighbours (no lingering `s="n"`
# style index) while the stored value/type stays a text string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.881.53"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").Value = "1.637.57"
$ws.Range("E3").Value = "  -1.06%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "215.31"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5031"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2559"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.58%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06378"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.66"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.54%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07727"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.75%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "4.256"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").Value = "1.642.01"
$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").Value = "1.864.47"
$ws.Range("E14").Value = "  -0.93%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.5447"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("D16").Value = "0.0₅7887"
$ws.Range("E16").Value = "  -1.70%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "64.24"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "25.872.68"
$ws.Range("E18").Value = "  -1.25%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "203.23"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -3.55%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.371"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.94%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "9.880"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.99%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.963"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.24%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.35%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.929"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +10.59%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "140.61"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1134"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "15.69"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.750"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.35%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.242"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("E31").Value = "  -3.08%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.268"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -3.13%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.178"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.52%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.540"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.62%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.368"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.73%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.632"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -3.85%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.8909"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.69%  "

$ws.Range("D38").Value = "1.158.72"
$ws.Range("E38").Value = "  -0.67%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.5607"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.12%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01566"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.42%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.47%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.639"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "99.67"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.67%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8044"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").Value = "1.776.66"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  +1.68%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4545"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.68%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "54.74"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.47%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.05056"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "

